$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- cryptos list refresh (GitHub Actions data update) ---
#
# Columns D (Price) and E (Volume(1h)) hold plain text (no custom number
# format on the cells). A handful of the refreshed Price values happen to
# parse cleanly as numbers (e.g. "1.00", "7.90"); a plain Value assignment
# would silently coerce those to a Number cell and drop the trailing
# zero(s). For just those cells we mark the cell as Text first so the
# literal string round-trips, then restore the original "Normal" style
# afterwards so no formatting change is left behind.

# Plain text / percent / link / coin-name cells -- safe to assign directly.
$ws.Range("D2").Value = "69.604.19"
$ws.Range("E2").Value = "  -4.03%  "
$ws.Range("D3").Value = "2.513.34"
$ws.Range("E3").Value = "  -5.60%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E5").Value = "  -2.45%  "
$ws.Range("E6").Value = "  -4.64%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -2.10%  "
$ws.Range("D9").Value = "2.513.39"
$ws.Range("E9").Value = "  -5.53%  "
$ws.Range("E10").Value = "  -6.28%  "
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("E12").Value = "  -5.30%  "
$ws.Range("E13").Value = "  -2.80%  "
$ws.Range("D14").Value = "2.970.28"
$ws.Range("E14").Value = "  -5.76%  "
$ws.Range("D15").Value = "69.469.38"
$ws.Range("E15").Value = "  -4.08%  "
$ws.Range("E16").Value = "  -4.61%  "
$ws.Range("E17").Value = "  -4.51%  "
$ws.Range("D18").Value = "2.513.91"
$ws.Range("E18").Value = "  -5.80%  "
$ws.Range("E19").Value = "  -7.34%  "
$ws.Range("E20").Value = "  -3.22%  "
$ws.Range("E21").Value = "  -5.26%  "
$ws.Range("E22").Value = "  -6.00%  "
$ws.Range("E23").Value = "  -4.02%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  -3.85%  "
$ws.Range("E26").Value = "  -6.11%  "
$ws.Range("E27").Value = "  -6.39%  "
$ws.Range("D28").Value = "2.640.48"
$ws.Range("E28").Value = "  -5.72%  "
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("D30").Value = "0.0₃0908"
$ws.Range("E30").Value = "  -5.11%  "
$ws.Range("E31").Value = "  -3.30%  "
$ws.Range("E32").Value = "  -2.91%  "
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("E34").Value = "  -3.49%  "
$ws.Range("E36").Value = "  -1.08%  "
$ws.Range("E37").Value = "  -3.59%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("E39").Value = "  -4.30%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("E42").Value = "  -3.10%  "
$ws.Range("E43").Value = "  -5.88%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("E44").Value = "  -8.35%  "
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E45").Value = "  -13.88%  "
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("E47").Value = "  -7.50%  "
$ws.Range("E48").Value = "  -4.39%  "
$ws.Range("E49").Value = "  -4.09%  "
$ws.Range("E50").Value = "  -5.36%  "
$ws.Range("E51").Value = "  -2.25%  "

# Numeric-looking Price cells: force Text format first so the exact
# string (incl. trailing zeros) is preserved instead of becoming a Number.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.65"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.66"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.512"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.340"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.85"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.98"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.46"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.77"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "353.11"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.95"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.98"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.13"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.986"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.90"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "484.45"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.30"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.116"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "155.98"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.90"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.62"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.78"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.320"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.36"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.18"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.26"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "143.17"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.54"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.529"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.61"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.589"

# Restore default styling on those cells -- only the values should differ from before.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"

